# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3
$ws.Range("C3").Value = 2400

# Row 4 - Intel(R) Dual Band Wireless-AC 8265 - 20.70.24.1
$ws.Range("C4").Value = 702
$ws.Range("D4").Value = 91.90000000000001

# Row 5 - Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7
$ws.Range("C5").Value = 520
$ws.Range("D5").Value = 93.7

# Row 6 - Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1
$ws.Range("C6").Value = 640
$ws.Range("D6").Value = 93.7

# Row 7 - Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3
$ws.Range("C7").Value = 278

# Row 8 and Row 9 swap the two "Dual Band Wireless-AC" driver entries
$ws.Range("A8").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.17.1"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 183

$ws.Range("A9").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.12.3"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 1131

# Row 10 - Intel(R) Dual Band Wireless-AC 8265 - 20.70.3.3
$ws.Range("C10").Value = 110

# Row 11 - Intel(R) Dual Band Wireless-AC 8265 - 20.50.3.3
$ws.Range("C11").Value = 669
$ws.Range("D11").Value = 98.3

# Row 12 - Totals
$ws.Range("C12").Value = 6633
